$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 30
$ws.Range("C19").Value = 0.01
$ws.Range("D19").Value = 0.003
$ws.Range("E19").Value = "Regular"
$ws.Range("F19").Value = "<function relu at 0x1104e69d8>"
$ws.Range("G19").Value = 0.8981999754905701
$ws.Range("H19").Value = 0.3424000144004822
$ws.Range("I19").Value = 0.04809999838471413
$ws.Range("J19").Value = 0.3424164950847626
$ws.Range("K19").Value = 3.330149173736572
$ws.Range("L19").Value = 0.3424000144004822
$ws.Range("M19").Value = "logs/results_121.log"
$ws.Range("N19").Value = "weights/model_121.ckpt"
$ws.Range("O19").Value = "tb/121/non_robust"
$ws.Range("P19").Value = "(1.4056362, 3.0716696, 6.664684, 11.19474, 13.153981, 10.241105, 6.051113)"
$ws.Range("Q19").Value = "(24.022184, 14.81685, 21.985281, 18.528912, 16.116236, 15.134631, 11.263296, 22.44767)"

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 30
$ws.Range("C20").Value = 0.02
$ws.Range("D20").Value = 0.0003
$ws.Range("E20").Value = "Regular"
$ws.Range("F20").Value = "<function relu at 0x10df1b9d8>"
$ws.Range("G20").Value = 0.9559000134468079
$ws.Range("H20").Value = 0.2694000005722046
$ws.Range("I20").Value = 0.05429999902844429
$ws.Range("J20").Value = 0.1658958792686462
$ws.Range("K20").Value = 4.698654651641846
$ws.Range("L20").Value = 0.2694000005722046
$ws.Range("M20").Value = "logs/results_123.log"
$ws.Range("N20").Value = "weights/model_123.ckpt"
$ws.Range("O20").Value = "tb/123/non_robust"
$ws.Range("P20").Value = "(0.3431952, 0.41664022, 0.5683927, 0.91492766, 1.3893795, 2.3776739, 3.8029735)"
$ws.Range("Q20").Value = "(6.3056436, 7.597117, 7.726599, 8.445029, 7.708117, 7.5270963, 7.9869733, 9.796146)"

# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 30
$ws.Range("C21").Value = 0.02
$ws.Range("D21").Value = 0.0003
$ws.Range("E21").Value = "Regular"
$ws.Range("F21").Value = "<function relu at 0x119b5f9d8>"
$ws.Range("G21").Value = 0.8652999997138977
$ws.Range("H21").Value = 0.1185000017285347
$ws.Range("I21").Value = 0.005900000222027302
$ws.Range("J21").Value = 0.4752624034881592
$ws.Range("K21").Value = 3.402863264083862
$ws.Range("L21").Value = 0.1185000017285347
$ws.Range("M21").Value = "logs/results_126.log"
$ws.Range("N21").Value = "weights/model_126.ckpt"
$ws.Range("O21").Value = "tb/126/non_robust"
$ws.Range("P21").Value = "(0.90816385, 0.7060737, 0.9453165, 1.4311132, 1.8074349, 2.626249, 3.3506863)"
$ws.Range("Q21").Value = "(15.063638, 2.4488645, 2.8450708, 3.4224455, 2.870352, 2.775072, 2.7803397, 3.6928034)"
Write-Host "done"
